$d = $word.ActiveDocument

# 1. Remove the existing _GoBack bookmark (currently located after
#    "Przybicie informacji do tablicy ogloszeniowej").
$existing = $d.Bookmarks("_GoBack")
$existing.Delete()

# 2. Locate the phrase that needs to be replaced, inside the paragraph
#    describing the exhibition manager's duties, and retype it.
$target = $d.Content
$found = $target.Find.Execute("zleca kierownikowi ochrony ustanowienie odpowiedniego grafiku na ustalone dni.")
$target.Text = "ochrony ustawia grafik ochrony na dane dni, które wystawa obejmuje."

# 3. The retyped text starts life as its own run, distinct from the
#    "Następnie kierownik " text that precedes it. Force that run break
#    with a temporary bookmark (removed again after a save commits the
#    split).
$afterKierownik = $d.Content
$afterKierownik.Find.Execute("kierownik ochrony ustawia grafik")
$splitPoint = $d.Range($afterKierownik.Start + 10, $afterKierownik.Start + 10)
$d.Bookmarks.Add("_TempSplit", $splitPoint)

# 4. Re-insert the _GoBack bookmark right after the newly typed text,
#    matching the place where the author's cursor ended up after the
#    edit.
$endOfEdit = $d.Content
$endOfEdit.Find.Execute("ochrony ustawia grafik ochrony na dane dni, które wystawa obejmuje.")
$mark = $d.Range($endOfEdit.End, $endOfEdit.End)
$d.Bookmarks.Add("_GoBack", $mark)

# 5. Commit the edits so far, then drop the temporary bookmark without
#    the surrounding runs collapsing back together.
$d.Save()
$tempBookmark = $d.Bookmarks("_TempSplit")
$tempBookmark.Delete()
